$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("J2").Value = 5.7
$ws.Range("K2").Value = 6.4
$ws.Range("L2").Value = 1.37
$ws.Range("Q2").Value = 1.79
$ws.Range("R2").Value = 1.42
$ws.Range("S2").Value = 3.1
$ws.Range("U2").Value = 1.7
$ws.Range("W2").Value = 3.75
$ws.Range("X2").Value = 22
$ws.Range("Z2").Value = 130
$ws.Range("AC2").Value = 16
$ws.Range("AD2").Value = 44
$ws.Range("AE2").Value = 260
$ws.Range("AF2").Value = 8.6
